$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3713.4285
$ws.Range("I86").Value = 3625
$ws.Range("J86").Value = 3831.3333
$ws.Range("K86").Value = 3625
$ws.Range("L86").Value = 3831.3333
$ws.Range("M86").Value = -2502
$ws.Range("N86").Value = -6077.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 3302.8462
$ws.Range("I88").Value = 8333.333000000001
$ws.Range("J88").Value = 1793.7
$ws.Range("K88").Value = 8333.333000000001
$ws.Range("L88").Value = 1793.7
$ws.Range("M88").Value = -7927.333000000001
$ws.Range("N88").Value = -2605.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3713.4285
$ws.Range("I89").Value = 3625
$ws.Range("J89").Value = 3831.3333
$ws.Range("K89").Value = 18125
$ws.Range("L89").Value = 19156.6665
$ws.Range("M89").Value = -12509
$ws.Range("N89").Value = -30388.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 3302.8462
$ws.Range("I91").Value = 8333.333000000001
$ws.Range("J91").Value = 1793.7
$ws.Range("K91").Value = 8333.333000000001
$ws.Range("L91").Value = 1793.7
$ws.Range("M91").Value = -6929.333000000001
$ws.Range("N91").Value = -4601.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1057.5106
$ws.Range("I98").Value = 1134.95
$ws.Range("K98").Value = 1134.95
$ws.Range("M98").Value = 363.05

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3474.25
$ws.Range("I113").Value = 2965.6667
$ws.Range("K113").Value = 2965.6667
$ws.Range("M113").Value = 288.3332999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 685.3889
$ws.Range("I115").Value = 716.75
$ws.Range("J115").Value = 434.5
$ws.Range("K115").Value = 2150.25
$ws.Range("L115").Value = 1303.5
$ws.Range("M115").Value = -583.25
$ws.Range("N115").Value = -4437.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1057.5106
$ws.Range("I122").Value = 1134.95
$ws.Range("K122").Value = 3404.85
$ws.Range("M122").Value = -954.8500000000004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 100021064
$ws.Range("I137").Value = 250000750
$ws.Range("J137").Value = 34601
$ws.Range("K137").Value = 750002250
$ws.Range("L137").Value = 103803
$ws.Range("M137").Value = -749999700
$ws.Range("N137").Value = -108903

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1052492.2
$ws.Range("I32").Value = 1236539.1
$ws.Range("K32").Value = 1236539.1
$ws.Range("M32").Value = -1236252.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2782.875
$ws.Range("I45").Value = 2217.3845
$ws.Range("K45").Value = 2217.3845
$ws.Range("M45").Value = -1840.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3452149
$ws.Range("I61").Value = 4048.8147
$ws.Range("J61").Value = 50001500
$ws.Range("K61").Value = 4048.8147
$ws.Range("L61").Value = 50001500
$ws.Range("M61").Value = -3836.8147
$ws.Range("N61").Value = -50001924

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2931.3584
$ws.Range("I132").Value = 2054.3547
$ws.Range("J132").Value = 4167.136
$ws.Range("K132").Value = 6163.0641
$ws.Range("L132").Value = 12501.408
$ws.Range("M132").Value = -3633.0641
$ws.Range("N132").Value = -17561.408

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 64614.4
$ws.Range("J133").Value = 64614.4
$ws.Range("L133").Value = 64614.4
$ws.Range("N133").Value = -69674.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3452149
$ws.Range("I136").Value = 4048.8147
$ws.Range("J136").Value = 50001500
$ws.Range("K136").Value = 12146.4441
$ws.Range("L136").Value = 150004500
$ws.Range("M136").Value = -9596.444100000001
$ws.Range("N136").Value = -150009600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 37500
$ws.Range("J28").Value = 37500
$ws.Range("L28").Value = 37500
$ws.Range("N28").Value = -38088

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 71042.25
$ws.Range("I58").Value = 39709
$ws.Range("J58").Value = 81486.664
$ws.Range("K58").Value = 39709
$ws.Range("L58").Value = 81486.664
$ws.Range("M58").Value = -39415
$ws.Range("N58").Value = -82074.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 78593.336
$ws.Range("J60").Value = 78593.336
$ws.Range("L60").Value = 78593.336
$ws.Range("N60").Value = -79791.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2226494
$ws.Range("I134").Value = 2505.1587
$ws.Range("K134").Value = 7515.4761
$ws.Range("M134").Value = -4980.4761

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 927825.7
$ws.Range("I31").Value = 1345921.1
$ws.Range("K31").Value = 1345921.1
$ws.Range("M31").Value = -1345626.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 927825.7
$ws.Range("I34").Value = 1345921.1
$ws.Range("K34").Value = 1345921.1
$ws.Range("M34").Value = -1345719.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 8433
$ws.Range("I122").Value = 1520.4615
$ws.Range("K122").Value = 4561.3845
$ws.Range("M122").Value = -2111.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2433.6
$ws.Range("I132").Value = 2133.7144
$ws.Range("K132").Value = 6401.1432
$ws.Range("M132").Value = -3871.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 69439.164
$ws.Range("I135").Value = 60000
$ws.Range("J135").Value = 69849.56
$ws.Range("K135").Value = 60000
$ws.Range("L135").Value = 69849.56
$ws.Range("M135").Value = -54930
$ws.Range("N135").Value = -79989.56

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 15153746
$ws.Range("I121").Value = 276.3684
$ws.Range("J121").Value = 35719172
$ws.Range("K121").Value = 829.1052
$ws.Range("L121").Value = 107157516
$ws.Range("M121").Value = 480.8948
$ws.Range("N121").Value = -107160136

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3089497.2
$ws.Range("I131").Value = 1242.8572
$ws.Range("J131").Value = 3834938
$ws.Range("K131").Value = 3728.5716
$ws.Range("L131").Value = 11504814
$ws.Range("M131").Value = 1311.4284
$ws.Range("N131").Value = -11514894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 43173.5
$ws.Range("J32").Value = 43173.5
$ws.Range("L32").Value = 43173.5
$ws.Range("N32").Value = -43765.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 40000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 40000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 40000
$ws.Range("N62").Value = -41372
$ws.Range("M62").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H65").Value = 40000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 40000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 120000
$ws.Range("N65").Value = -126864
$ws.Range("M65").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1833.24
$ws.Range("I102").Value = 1942.5264
$ws.Range("K102").Value = 1942.5264
$ws.Range("M102").Value = -320.5264

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3231.9556
$ws.Range("J122").Value = 3106.8
$ws.Range("L122").Value = 9320.400000000001
$ws.Range("N122").Value = -14220.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 10074.038
$ws.Range("I132").Value = 10567.708
$ws.Range("J132").Value = 4150
$ws.Range("K132").Value = 31703.124
$ws.Range("L132").Value = 12450
$ws.Range("M132").Value = -29173.124
$ws.Range("N132").Value = -17510

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2136.4285
$ws.Range("I40").Value = 1952.8
$ws.Range("J40").Value = 3666.6667
$ws.Range("K40").Value = 1952.8
$ws.Range("L40").Value = 3666.6667
$ws.Range("M40").Value = -1816.8
$ws.Range("N40").Value = -3938.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 47667.445
$ws.Range("J74").Value = 53769
$ws.Range("L74").Value = 53769
$ws.Range("N74").Value = -55765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H77").Value = 47667.445
$ws.Range("J77").Value = 53769
$ws.Range("L77").Value = 161307
$ws.Range("N77").Value = -171291

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3888.3235
$ws.Range("I122").Value = 2825.8572
$ws.Range("J122").Value = 5604.615
$ws.Range("K122").Value = 8477.571599999999
$ws.Range("L122").Value = 16813.845
$ws.Range("M122").Value = -6027.571599999999
$ws.Range("N122").Value = -21713.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 31795.256
$ws.Range("I122").Value = 2845.1562
$ws.Range("K122").Value = 8535.4686
$ws.Range("M122").Value = -6085.4686
